# Daily attendance processing - 2026-01-25 17:58:11
# Normalize the "Recorded By" column (G) so that entries listing both the
# System and the user's email are ordered as "email, System" instead of
# "System, email".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
